$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold values that look like plain numbers (e.g. "0.630",
# "22.60"), but the source data stores them as literal text. Pre-format them
# as Text so Excel does not silently convert them to numbers (which would
# drop significant trailing zeros) when we assign the new values below.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D13", "D15", "D16", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.164.78'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = '1.999.23'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '246.53'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").Value = '0.630'
$ws.Range("E6").Value = '  +2.86%  '
$ws.Range("D7").Value = '60.21'
$ws.Range("E7").Value = '  +4.23%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").Value = '0.0800'
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '14.94'
$ws.Range("E12").Value = '  +8.10%  '
$ws.Range("D13").Value = '22.60'
$ws.Range("E13").Value = '  +7.36%  '
$ws.Range("D14").Value = '2.296.79'
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").Value = '0.845'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '5.45'
$ws.Range("E16").Value = '  +3.15%  '
$ws.Range("D17").Value = '2.016.49'
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").Value = '37.096.49'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("D19").Value = '70.29'
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").Value = '0.0₃0862'
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("D21").Value = '5.17'
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").Value = '230.36'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '2.47'
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").Value = '9.38'
$ws.Range("E26").Value = '  +2.97%  '
$ws.Range("E27").Value = '  +5.69%  '
$ws.Range("D28").Value = '163.26'
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").Value = '19.65'
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("E30").Value = '  +14.12%  '
$ws.Range("E31").Value = '  +1.36%  '
$ws.Range("D32").Value = '4.81'
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").Value = '0.0658'
$ws.Range("E33").Value = '  +8.31%  '
$ws.Range("D34").Value = '4.50'
$ws.Range("E34").Value = '  +3.08%  '
$ws.Range("D35").Value = '2.37'
$ws.Range("E35").Value = '  +5.70%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '1.81'
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("D38").Value = '3.26'
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").Value = '5.41'
$ws.Range("E39").Value = '  +4.02%  '
$ws.Range("D40").Value = '0.0978'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("E43").Value = '  +2.25%  '
$ws.Range("D44").Value = '16.66'
$ws.Range("E44").Value = '  +5.86%  '
$ws.Range("D45").Value = '90.89'
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("D46").Value = '1.376.00'
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("D47").Value = '1.04'
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("D48").Value = '7.26'
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '2.03'
$ws.Range("E50").Value = '  +16.13%  '
$ws.Range("D51").Value = '46.24'
$ws.Range("E51").Value = '  +5.58%  '

# Drop the temporary Text number-format override so the cell styling
# matches the original (unstyled) data cells again.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}

